$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iterations")

$ws.Range("B2").Value = 112186436.73
$ws.Range("D2").Value = 949876.27
$ws.Range("E2").Value = 9049342.050000001
$ws.Range("F2").Value = -273433.18
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 12243394.18
$ws.Range("I2").Value = 2646153.86
$ws.Range("J2").Value = 59299813.38
$ws.Range("K2").Value = 8140315.9
$ws.Range("L2").Value = 19394703.97
$ws.Range("M2").Value = 15751661.2
$ws.Range("N2").Value = 4078063.85
$ws.Range("O2").Value = 20078585.44
$ws.Range("Q2").Value = 199410823.28
$ws.Range("R2").Value = 33394142.52
$ws.Range("S2").Value = 96152506.36
$ws.Range("T2").Value = 12055622.7
$ws.Range("U2").Value = 14067054.58
$ws.Range("V2").Value = 35458458.55
$ws.Range("W2").Value = 15062319.49
$ws.Range("X2").Value = 58919491.72
$ws.Range("Y2").Value = -155477.56
$ws.Range("Z2").Value = 2941036.87
$ws.Range("AA2").Value = 5240354.25
$ws.Range("AC2").Value = 25257955.68
$ws.Range("AD2").Value = 30408914.89
$ws.Range("AE2").Value = 17939366.72
$ws.Range("AF2").Value = 6009155.18
$ws.Range("AG2").Value = 6873271.76
$ws.Range("AH2").Value = 66746670.95
$ws.Range("AI2").Value = 31876486.84
$ws.Range("AJ2").Value = 228402877.46
$ws.Range("AK2").Value = 0
$ws.Range("AL2").Value = 31432003.97
$ws.Range("AM2").Value = 17724790
$ws.Range("AO2").Value = 6459880.81
$ws.Range("AP2").Value = 472425.44
$ws.Range("AQ2").Value = 54823649.7
$ws.Range("AR2").Value = 10159480.5
$ws.Range("AS2").Value = 45480575.49
$ws.Range("AT2").Value = 10027318.81
$ws.Range("AU2").Value = 1281963.45
$ws.Range("AW2").Value = -161935.12
$ws.Range("AX2").Value = 10902861.83
$ws.Range("AY2").Value = 1195942.67
$ws.Range("AZ2").Value = 59851334.36
$ws.Range("BA2").Value = 3620250.88
$ws.Range("BB2").Value = 30782005.6
$ws.Range("BC2").Value = 62052071.28
$ws.Range("BD2").Value = -350217.95
$ws.Range("BE2").Value = 98348347.98999999
$ws.Range("BG2").Value = 2272433.25
$ws.Range("BH2").Value = 516989963.38
$ws.Range("BJ2").Value = 61306732.85
$ws.Range("BK2").Value = 8646051.92
$ws.Range("BL2").Value = 58907525.44
$ws.Range("BM2").Value = 102624872.94
$ws.Range("BN2").Value = 14432633.82
$ws.Range("BO2").Value = 7337269.6
$ws.Range("BP2").Value = 11793094.31
$ws.Range("BQ2").Value = 235247.08
$ws.Range("BR2").Value = 107563381.61
$ws.Range("BS2").Value = 3673430.39
$ws.Range("BT2").Value = 112656242.18
$ws.Range("BV2").Value = 23690234.13
$ws.Range("BW2").Value = 10785328.76
$ws.Range("BX2").Value = 78521210.66
$ws.Range("CA2").Value = 15522899.32
$ws.Range("CB2").Value = 0
$ws.Range("CC2").Value = 0
$ws.Range("CD2").Value = 36759931.57
$ws.Range("CE2").Value = 10372366.58
$ws.Range("CF2").Value = 50757641.24
$ws.Range("CG2").Value = 2195049.46
$ws.Range("CJ2").Value = 166766207.22
$ws.Range("CK2").Value = 32974031.03
$ws.Range("CL2").Value = 12583245.94
$ws.Range("CM2").Value = 2471782.27
$ws.Range("CN2").Value = 25451101.36
$ws.Range("CO2").Value = 18793746.08
$ws.Range("CP2").Value = 61075352.18
$ws.Range("CQ2").Value = 4873958.19
$ws.Range("CR2").Value = 10784133.13
$ws.Range("CS2").Value = 34619060.22
$ws.Range("CT2").Value = 214352530.27
$ws.Range("CU2").Value = 4590835.49
$ws.Range("CV2").Value = -120114.76
$ws.Range("CW2").Value = 10444527.3
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = -25148.61
$ws.Range("D3").Value = 437018.85
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 76109.34
$ws.Range("G3").Value = -357020.79
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("N3").Value = 258232.23
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -778.11
$ws.Range("R3").Value = 71987.83
$ws.Range("S3").Value = 104886.84
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 659659.23
$ws.Range("W3").Value = 545860.3199999999
$ws.Range("X3").Value = 1223019.43
$ws.Range("Y3").Value = 140874.33
$ws.Range("Z3").Value = 269024.24
$ws.Range("AA3").Value = 241593.45
$ws.Range("AB3").Value = -25534.63
$ws.Range("AC3").Value = 836263.5600000001
$ws.Range("AD3").Value = 457064.27
$ws.Range("AF3").Value = 732086.1800000001
$ws.Range("AG3").Value = 167352.85
$ws.Range("AI3").Value = 352497.87
$ws.Range("AJ3").Value = 141506.21
$ws.Range("AK3").Value = -16451.3
$ws.Range("AL3").Value = 156438.46
$ws.Range("AM3").Value = 122986.8
$ws.Range("AN3").Value = -91369.7
$ws.Range("AO3").Value = 0
$ws.Range("AP3").Value = 74434.60000000001
$ws.Range("AQ3").Value = 19402.44
$ws.Range("AR3").Value = 395686.59
$ws.Range("AT3").Value = 299026.82
$ws.Range("AU3").Value = 0
$ws.Range("AV3").Value = -57927.76
$ws.Range("AW3").Value = 33192.07
$ws.Range("AX3").Value = 202866.61
$ws.Range("AY3").Value = 43725.19
$ws.Range("AZ3").Value = 488256.3
$ws.Range("BA3").Value = 0
$ws.Range("BB3").Value = 655927.73
$ws.Range("BC3").Value = 235643.65
$ws.Range("BD3").Value = 215285.88
$ws.Range("BF3").Value = -5146.65
$ws.Range("BG3").Value = 0
$ws.Range("BH3").Value = 869853.0600000001
$ws.Range("BI3").Value = 18088.14
$ws.Range("BJ3").Value = 211571.96
$ws.Range("BK3").Value = 419671.38
$ws.Range("BM3").Value = 0
$ws.Range("BN3").Value = 591690.14
$ws.Range("BO3").Value = 78435.10000000001
$ws.Range("BP3").Value = 89068.63
$ws.Range("BR3").Value = 767001.02
$ws.Range("BT3").Value = 583882.85
$ws.Range("BU3").Value = -83860.10000000001
$ws.Range("BV3").Value = 52739.96
$ws.Range("BW3").Value = 82880.78999999999
$ws.Range("BX3").Value = 54160.26
$ws.Range("BY3").Value = -137207.23
$ws.Range("BZ3").Value = -226583.76
$ws.Range("CA3").Value = 253643.67
$ws.Range("CB3").Value = -1667.35
$ws.Range("CC3").Value = 355324.13
$ws.Range("CD3").Value = 0
$ws.Range("CE3").Value = 0
$ws.Range("CF3").Value = 0
$ws.Range("CG3").Value = 244312.17
$ws.Range("CH3").Value = -128070.79
$ws.Range("CI3").Value = -124743.69
$ws.Range("CJ3").Value = 380608.53
$ws.Range("CK3").Value = 1123450.63
$ws.Range("CL3").Value = 0
$ws.Range("CM3").Value = 48790.99
$ws.Range("CN3").Value = 687033.67
$ws.Range("CO3").Value = 258414.87
$ws.Range("CP3").Value = 0
$ws.Range("CQ3").Value = 0
$ws.Range("CR3").Value = 456527.35
$ws.Range("CS3").Value = 307667.13
$ws.Range("CU3").Value = 213558.1
$ws.Range("CW3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("AB4").Value = -40967.78
$ws.Range("AH4").Value = 0
$ws.Range("AX4").Value = 0
$ws.Range("BI4").Value = -23644.44
$ws.Range("BW4").Value = 0
$ws.Range("B6").Value = 112186436.73
$ws.Range("C6").Value = -25148.61
$ws.Range("D6").Value = 1386895.11
$ws.Range("E6").Value = 9049342.050000001
$ws.Range("F6").Value = -197323.84
$ws.Range("G6").Value = -357020.79
$ws.Range("H6").Value = 12243394.18
$ws.Range("I6").Value = 2646153.86
$ws.Range("J6").Value = 59299813.38
$ws.Range("K6").Value = 8140315.9
$ws.Range("L6").Value = 19394703.97
$ws.Range("M6").Value = 15751661.2
$ws.Range("N6").Value = 4336296.08
$ws.Range("O6").Value = 20078585.44
$ws.Range("P6").Value = -778.11
$ws.Range("Q6").Value = 199410823.28
$ws.Range("R6").Value = 33466130.35
$ws.Range("S6").Value = 96257393.20999999
$ws.Range("T6").Value = 12055622.7
$ws.Range("U6").Value = 14067054.58
$ws.Range("V6").Value = 36118117.78
$ws.Range("W6").Value = 15608179.81
$ws.Range("X6").Value = 60142511.15
$ws.Range("Y6").Value = -14603.24
$ws.Range("Z6").Value = 3210061.11
$ws.Range("AA6").Value = 5481947.7
$ws.Range("AB6").Value = -66502.39999999999
$ws.Range("AC6").Value = 26094219.23
$ws.Range("AD6").Value = 30865979.16
$ws.Range("AE6").Value = 17939366.72
$ws.Range("AF6").Value = 6741241.36
$ws.Range("AG6").Value = 7040624.62
$ws.Range("AH6").Value = 66746670.95
$ws.Range("AI6").Value = 32228984.71
$ws.Range("AJ6").Value = 228544383.68
$ws.Range("AK6").Value = -16451.3
$ws.Range("AL6").Value = 31588442.43
$ws.Range("AM6").Value = 17847776.8
$ws.Range("AN6").Value = -91369.7
$ws.Range("AO6").Value = 6459880.81
$ws.Range("AP6").Value = 546860.03
$ws.Range("AQ6").Value = 54843052.13
$ws.Range("AR6").Value = 10555167.1
$ws.Range("AS6").Value = 45480575.49
$ws.Range("AT6").Value = 10326345.63
$ws.Range("AU6").Value = 1281963.45
$ws.Range("AV6").Value = -57927.76
$ws.Range("AW6").Value = -128743.05
$ws.Range("AX6").Value = 11105728.43
$ws.Range("AY6").Value = 1239667.86
$ws.Range("AZ6").Value = 60339590.66
$ws.Range("BA6").Value = 3620250.88
$ws.Range("BB6").Value = 31437933.33
$ws.Range("BC6").Value = 62287714.93
$ws.Range("BD6").Value = -134932.07
$ws.Range("BE6").Value = 98348347.98999999
$ws.Range("BF6").Value = -5146.65
$ws.Range("BG6").Value = 2272433.25
$ws.Range("BH6").Value = 517859816.43
$ws.Range("BI6").Value = -5556.3
$ws.Range("BJ6").Value = 61518304.81
$ws.Range("BK6").Value = 9065723.300000001
$ws.Range("BL6").Value = 58907525.44
$ws.Range("BM6").Value = 102624872.94
$ws.Range("BN6").Value = 15024323.96
$ws.Range("BO6").Value = 7415704.7
$ws.Range("BP6").Value = 11882162.94
$ws.Range("BQ6").Value = 235247.08
$ws.Range("BR6").Value = 108330382.63
$ws.Range("BS6").Value = 3673430.39
$ws.Range("BT6").Value = 113240125.03
$ws.Range("BU6").Value = -83860.10000000001
$ws.Range("BV6").Value = 23742974.1
$ws.Range("BW6").Value = 10868209.55
$ws.Range("BX6").Value = 78575370.92
$ws.Range("BY6").Value = -137207.23
$ws.Range("BZ6").Value = -226583.76
$ws.Range("CA6").Value = 15776542.99
$ws.Range("CB6").Value = -1667.35
$ws.Range("CC6").Value = 355324.13
$ws.Range("CD6").Value = 36759931.57
$ws.Range("CE6").Value = 10372366.58
$ws.Range("CF6").Value = 50757641.24
$ws.Range("CG6").Value = 2439361.63
$ws.Range("CH6").Value = -128070.79
$ws.Range("CI6").Value = -124743.69
$ws.Range("CJ6").Value = 167146815.75
$ws.Range("CK6").Value = 34097481.66
$ws.Range("CL6").Value = 12583245.94
$ws.Range("CM6").Value = 2520573.26
$ws.Range("CN6").Value = 26138135.03
$ws.Range("CO6").Value = 19052160.95
$ws.Range("CP6").Value = 61075352.18
$ws.Range("CQ6").Value = 4873958.19
$ws.Range("CR6").Value = 11240660.48
$ws.Range("CS6").Value = 34926727.35
$ws.Range("CT6").Value = 214352530.27
$ws.Range("CU6").Value = 4804393.59
$ws.Range("CV6").Value = -120114.76
$ws.Range("CW6").Value = 10444527.3
